# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the handback
# files for 1232741e-... and 362a2034-... are ready for handoff again,
# refreshes the relevant timestamps, records "version not latest" error
# details for both locales, and widens the Error Detail column.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

$msg1232741e = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9444296fab8b7f248322409b01a7bf5374342db/e2e/1232741e-5f8e-43f4-be1c-8ecbae8666f7.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0904d5321658a866ae61da0d11b42906a980edab/e2e/1232741e-5f8e-43f4-be1c-8ecbae8666f7.md."
$msg362a2034 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9444296fab8b7f248322409b01a7bf5374342db/e2e/362a2034-6cc5-4300-be4f-b799ccdc07c4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0904d5321658a866ae61da0d11b42906a980edab/e2e/362a2034-6cc5-4300-be4f-b799ccdc07c4.md."

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $status
$wsOverview.Range("F4").Value = $status
$wsOverview.Range("G4").Value = "2016-08-30 00:28:15"

$wsOverview.Range("E5").Value = $status
$wsOverview.Range("F5").Value = $status
$wsOverview.Range("G5").Value = "2016-08-30 00:28:15"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C4").Value = $status
$wsZhCn.Range("H4").Value = "2016-08-30 00:28:11"
$wsZhCn.Range("P4").Value = $msg1232741e

$wsZhCn.Range("C5").Value = $status
$wsZhCn.Range("H5").Value = "2016-08-30 00:28:11"
$wsZhCn.Range("P5").Value = $msg362a2034

$wsZhCn.Columns.Item(16).ColumnWidth = 39.142857142857146

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C4").Value = $status
$wsDeDe.Range("H4").Value = "2016-08-30 00:28:15"
$wsDeDe.Range("P4").Value = $msg1232741e

$wsDeDe.Range("C5").Value = $status
$wsDeDe.Range("H5").Value = "2016-08-30 00:28:15"
$wsDeDe.Range("P5").Value = $msg362a2034

$wsDeDe.Columns.Item(16).ColumnWidth = 39.142857142857146
